# Generate Report for Handback
# This script applies the "handback" status update to the localization-status
# workbook: marks zh-cn/de-de as handed back, fills in the Latest Target File /
# Latest Handback File / Latest Handback DateTime columns for both language
# sheets, widens a few columns to fit the new content, and adds hyperlinks
# for the newly populated "Latest Target File" cells.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Hyperlink colour used by the workbook's built-in "HyperLink" cell style
# (font color FF6495ED, single underline) -- OLE color is 0xBBGGRR reversed,
# i.e. just the RRGGBB value read as an integer.
$hyperlinkColor = 0x6495ED

function Apply-HyperlinkLook($range) {
    $range.Font.Underline = 2   # xlUnderlineStyleSingle
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns and update the
# status text shown there (shared with the per-language sheets' Status column)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# "Status" column shares the same underlying text as the Overview sheet's
# zh-cn/de-de columns ("In Translation" -> "Handed back: in sync with en-US")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$zhMdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e159470bc911f257b150f57ea163b2fd7fbc4490/e2e/41fe5f37-68fd-4cb3-9f36-352333e89f54.md"
$zhMdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e159470bc911f257b150f57ea163b2fd7fbc4490/e2e/b7b6e70a-dae7-45b7-98c2-0f6375d9417b.md"

$wsZh.Range("I2").Value = "41fe5f37-68fd-4cb3-9f36-352333e89f54.md"
Apply-HyperlinkLook $wsZh.Range("I2")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhMdUrl1, "", "", "41fe5f37-68fd-4cb3-9f36-352333e89f54.md") | Out-Null
$wsZh.Range("J2").Value = "41fe5f37-68fd-4cb3-9f36-352333e89f54.ce1bd6e4eded521268699596770927979c56b5a9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-22 14:25:40"

$wsZh.Range("I3").Value = "b7b6e70a-dae7-45b7-98c2-0f6375d9417b.md"
Apply-HyperlinkLook $wsZh.Range("I3")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhMdUrl2, "", "", "b7b6e70a-dae7-45b7-98c2-0f6375d9417b.md") | Out-Null
$wsZh.Range("J3").Value = "b7b6e70a-dae7-45b7-98c2-0f6375d9417b.d380fdfbf37b19f918608cb8269567a0c76f4e6b.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-22 14:25:40"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$deMdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e159470bc911f257b150f57ea163b2fd7fbc4490/e2e/41fe5f37-68fd-4cb3-9f36-352333e89f54.md"
$deMdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e159470bc911f257b150f57ea163b2fd7fbc4490/e2e/b7b6e70a-dae7-45b7-98c2-0f6375d9417b.md"

$wsDe.Range("I2").Value = "41fe5f37-68fd-4cb3-9f36-352333e89f54.md"
Apply-HyperlinkLook $wsDe.Range("I2")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deMdUrl1, "", "", "41fe5f37-68fd-4cb3-9f36-352333e89f54.md") | Out-Null
$wsDe.Range("J2").Value = "41fe5f37-68fd-4cb3-9f36-352333e89f54.ce1bd6e4eded521268699596770927979c56b5a9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-22 14:25:47"

$wsDe.Range("I3").Value = "b7b6e70a-dae7-45b7-98c2-0f6375d9417b.md"
Apply-HyperlinkLook $wsDe.Range("I3")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deMdUrl2, "", "", "b7b6e70a-dae7-45b7-98c2-0f6375d9417b.md") | Out-Null
$wsDe.Range("J3").Value = "b7b6e70a-dae7-45b7-98c2-0f6375d9417b.d380fdfbf37b19f918608cb8269567a0c76f4e6b.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-22 14:25:47"
